# Auto-generated edit script
# Refreshes the crypto price table: Price (D), Volume(1h)% (E), and
# Hora (G) columns for data rows 2-51, per the symbol-list update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces Excel to store the value as text,
# matching the workbook's existing text-typed numeric-looking cells.
$ws.Range("D2").Value = "'321.72"
$ws.Range("E2").Value = "'8.00%"
$ws.Range("G2").Value = "'9"
$ws.Range("D3").Value = "'48.58"
$ws.Range("E3").Value = "'15.49%"
$ws.Range("G3").Value = "'9"
$ws.Range("D4").Value = "'5.230"
$ws.Range("E4").Value = "'4.37%"
$ws.Range("G4").Value = "'9"
$ws.Range("D5").Value = "'0.08094"
$ws.Range("E5").Value = "'7.59%"
$ws.Range("G5").Value = "'9"
$ws.Range("D6").Value = "'4.595"
$ws.Range("E6").Value = "'5.20%"
$ws.Range("G6").Value = "'9"
$ws.Range("D7").Value = "'1.639"
$ws.Range("E7").Value = "'2.73%"
$ws.Range("G7").Value = "'9"
$ws.Range("D8").Value = "'1.201"
$ws.Range("E8").Value = "'30.78%"
$ws.Range("G8").Value = "'9"
$ws.Range("D9").Value = "'0.1297"
$ws.Range("E9").Value = "'9.42%"
$ws.Range("G9").Value = "'9"
$ws.Range("D10").Value = "'0.1951"
$ws.Range("E10").Value = "'6.57%"
$ws.Range("G10").Value = "'9"
$ws.Range("D11").Value = "'0.09465"
$ws.Range("E11").Value = "'5.90%"
$ws.Range("G11").Value = "'9"
$ws.Range("D12").Value = "'0.04644"
$ws.Range("E12").Value = "'12.48%"
$ws.Range("G12").Value = "'9"
$ws.Range("D13").Value = "'0.1049"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("G13").Value = "'9"
$ws.Range("D14").Value = "'0.001319"
$ws.Range("E14").Value = "'3.26%"
$ws.Range("G14").Value = "'9"
$ws.Range("D15").Value = "'0.005883"
$ws.Range("E15").Value = "'-2.11%"
$ws.Range("G15").Value = "'9"
$ws.Range("D16").Value = "'3.341"
$ws.Range("E16").Value = "'-0.07%"
$ws.Range("G16").Value = "'9"
$ws.Range("D17").Value = "'2.429"
$ws.Range("E17").Value = "'1.14%"
$ws.Range("G17").Value = "'9"
$ws.Range("D18").Value = "'0.3402"
$ws.Range("E18").Value = "'2.11%"
$ws.Range("G18").Value = "'9"
$ws.Range("D19").Value = "'8.128"
$ws.Range("E19").Value = "'-2.88%"
$ws.Range("G19").Value = "'9"
$ws.Range("D20").Value = "'0.1411"
$ws.Range("E20").Value = "'4.34%"
$ws.Range("G20").Value = "'9"
$ws.Range("E21").Value = "'0.69%"
$ws.Range("G21").Value = "'9"
$ws.Range("E22").Value = "'4.86%"
$ws.Range("G22").Value = "'9"
$ws.Range("D23").Value = "'0.001307"
$ws.Range("E23").Value = "'3.17%"
$ws.Range("G23").Value = "'9"
$ws.Range("D24").Value = "'0.004249"
$ws.Range("E24").Value = "'9.19%"
$ws.Range("G24").Value = "'9"
$ws.Range("E25").Value = "'3.73%"
$ws.Range("G25").Value = "'9"
$ws.Range("D26").Value = "'0.0003543"
$ws.Range("E26").Value = "'-4.86%"
$ws.Range("G26").Value = "'9"
$ws.Range("G27").Value = "'9"
$ws.Range("G28").Value = "'9"
$ws.Range("G29").Value = "'9"
$ws.Range("G30").Value = "'9"
$ws.Range("G31").Value = "'9"
$ws.Range("G32").Value = "'9"
$ws.Range("G33").Value = "'9"
$ws.Range("G34").Value = "'9"
$ws.Range("G35").Value = "'9"
$ws.Range("G36").Value = "'9"
$ws.Range("G37").Value = "'9"
$ws.Range("D38").Value = "'0.02657"
$ws.Range("E38").Value = "'11.25%"
$ws.Range("G38").Value = "'9"
$ws.Range("D39").Value = "'0.05620"
$ws.Range("E39").Value = "'7.47%"
$ws.Range("G39").Value = "'9"
$ws.Range("D40").Value = "'0.006307"
$ws.Range("E40").Value = "'-10.44%"
$ws.Range("G40").Value = "'9"
$ws.Range("D41").Value = "'0.007694"
$ws.Range("E41").Value = "'-0.84%"
$ws.Range("G41").Value = "'9"
$ws.Range("D42").Value = "'0.1437"
$ws.Range("E42").Value = "'8.31%"
$ws.Range("G42").Value = "'9"
$ws.Range("D43").Value = "'0.007700"
$ws.Range("E43").Value = "'3.93%"
$ws.Range("G43").Value = "'9"
$ws.Range("E44").Value = "'13.43%"
$ws.Range("G44").Value = "'9"
$ws.Range("E45").Value = "'-0.91%"
$ws.Range("G45").Value = "'9"
$ws.Range("D46").Value = "'0.00007016"
$ws.Range("E46").Value = "'6.87%"
$ws.Range("G46").Value = "'9"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("G47").Value = "'9"
$ws.Range("D48").Value = "'0.05667"
$ws.Range("E48").Value = "'24.83%"
$ws.Range("G48").Value = "'9"
$ws.Range("D49").Value = "'0.004004"
$ws.Range("E49").Value = "'-4.81%"
$ws.Range("G49").Value = "'9"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("G50").Value = "'9"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.05%"
$ws.Range("G51").Value = "'9"
